$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Beta_M1_RN"
$ws.Range("C1").Value = "Beta_CM2_RN"
$ws.Range("D1").Value = "Beta_CMN3_RN"
$ws.Range("E1").Value = "Beta_CMN4_RN"
